$d = $word.ActiveDocument

# Locate the paragraph that ends the "KODLAR" section with the Beulian sprite-change
# fix ("...sprite degisimi duzeltilecek."), using an accent-insensitive / ASCII substring
# so the search is robust regardless of exact Unicode normalization.
$target = $d.Content.Duplicate
[void]$target.Find.Execute("sprite", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$target.Expand(4)

$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $target.Start) {
        $anchorIndex = $i
        break
    }
}

# Insert a brand-new paragraph right after it.
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$newParaIndex = $anchorIndex + 1
$insertPos = $d.Paragraphs.Item($newParaIndex).Range.Start

# Run 1
$run1 = $d.Range($insertPos, $insertPos)
$run1.InsertAfter("Minyoo timing platformların üstündeyken ")

# Run 2
$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("platform aniden yok olursa minyoo ")
$bookmarkPos = $run2.End

# Insert run 3's text before placing the bookmark: in this engine, inserting text exactly at a
# bookmark's (zero-width) position pushes the bookmark to the end of the new text, so the
# bookmark must be (re)created only after all surrounding text already exists.
$run3 = $d.Range($bookmarkPos, $bookmarkPos)
$run3.InsertAfter("bug’a girmesi sorunu çözülecek.")

# The "_GoBack" bookmark must move here from its old location (end of the "Pixo'nun..."
# paragraph further down); delete the old one first since bookmark names are unique.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
